# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from serial 45172 (2023-09-03) to serial 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45172) {
        $cell.Value2 = 45175
    }
}
